$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$vals = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)
for ($i=1; $i -le 12; $i++) {
    $c = $tcs.Colors($i)
    $c.RGB = $vals[$i-1]
}
Write-Output "done"
